$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: Status columns for zh-cn (B2) and de-de (C2), and Latest Handoff Date (D2)
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-30-19 14:30:56"

# zh-cn sheet: Status (C2) and Latest Handoff Datetime (E2)
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("E2").Value = "2016-03-19 14:30:53"

# de-de sheet: Status (C2) and Latest Handoff Datetime (E2)
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("E2").Value = "2016-03-19 14:30:56"
